# Remove specific rows from the "Test" sheet of the feature-space workbook.
# Rows being removed (1-based, original layout):
#   3  - download.jpg                       (Unknown - Test split)
#   4  - idi-amin-u4.jpg                    (Unknown - Test split)
#   12 - blank separator row
#   17 - IMG_20230305_175457.jpg            (Sameed - Test split)
#   20 - IMG_20230305_175508.jpg            (Sameed - Test split)
#   23 - blank separator row
# Deleting them (entire row, shifting cells up) compacts the table from
# A1:D31 down to A1:D25, matching the target layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(23, 20, 17, 12, 4, 3)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

$ws.Range("E14").Select()
